$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 2023 September
$ws.Range("B2").Value = 64
$ws.Range("C2").Value = 2

# Row 3 - 2023 October
$ws.Range("C3").Value = 9

# Row 5 - 2023 December
$ws.Range("B5").Value = 51
$ws.Range("C5").Value = 9

# Row 6 - 2024 January
$ws.Range("B6").Value = 52
$ws.Range("C6").Value = 9

# Row 8 - 2024 March
$ws.Range("B8").Value = 31
$ws.Range("C8").Value = 12

# Row 9 - 2024 April
$ws.Range("B9").Value = 54
$ws.Range("C9").Value = 13

# Row 10 - 2024 May
$ws.Range("B10").Value = 37
$ws.Range("C10").Value = 24

# Row 11 - 2024 June
$ws.Range("B11").Value = 46
$ws.Range("C11").Value = 25

# Row 12 - 2024 July
$ws.Range("B12").Value = 33
$ws.Range("C12").Value = 20

# Row 13 - 2024 August
$ws.Range("B13").Value = 14
$ws.Range("C13").Value = 19
